{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph that precedes it) that the site generator used\n// to append after the course's \"Requisitos\" list.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst JUPITER_TEXT = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst COPYRIGHT_SNIPPET = \"Powered by Jekyll and Github pages\";\n\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (jupiterIndex === -1 && text === JUPITER_TEXT) {\n    jupiterIndex = i;\n  }\n  if (copyrightIndex === -1 && text.indexOf(COPYRIGHT_SNIPPET) !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (jupiterIndex === -1 || copyrightIndex === -1) {\n  throw new Error(\"Could not locate the 'Ver no Jupiter' / copyright paragraphs to remove.\");\n}\n\n// The blank paragraph immediately preceding the \"Ver no Jupiter\" line is\n// part of the same trailing block and is removed along with it.\nconst blankIndex = jupiterIndex - 1;\n\n// Delete from the bottom up so earlier indices stay valid.\nconst indicesToDelete = [copyrightIndex, jupiterIndex, blankIndex].sort((a, b) => b - a);\nfor (const idx of indicesToDelete) {\n  paragraphs.items[idx].delete();\n}\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"(c) 2020 ...\" footer block\n# (and the blank paragraph that precedes it) that the site generator used\n# to append after the course's \"Requisitos\" list.\n$d = $word.ActiveDocument\n\n$jupiterIndex = -1\n$copyrightIndex = -1\n\n$total = $d.Paragraphs.Count\nfor ($i = 1; $i -le $total; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    if ($jupiterIndex -eq -1 -and $text -like \"*Ver no Jupiter Salvar em pdf Salvar em docx*\") {\n        $jupiterIndex = $i\n    }\n    if ($copyrightIndex -eq -1 -and $text -like \"*Powered by Jekyll and Github pages*\") {\n        $copyrightIndex = $i\n    }\n}\n\nif ($jupiterIndex -eq -1 -or $copyrightIndex -eq -1) {\n    throw \"Could not locate the 'Ver no Jupiter' / copyright paragraphs to remove.\"\n}\n\n# The blank paragraph immediately preceding the \"Ver no Jupiter\" line is\n# part of the same trailing block and is removed along with it.\n$blankIndex = $jupiterIndex - 1\n\n# Delete from the bottom up so earlier indices stay valid.\n$d.Paragraphs.Item($copyrightIndex).Range.Delete()\n$d.Paragraphs.Item($jupiterIndex).Range.Delete()\n$d.Paragraphs.Item($blankIndex).Range.Delete()\n"}
